$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2: update Authors value (cleaned data - extra whitespace normalization)
$ws.Range("E2").Value = '[Atas%Jenny%coreGivesNoEmail%1,                Bandy%Kenneth%coreGivesNoEmail%1,                Bradin%Stuart A.%coreGivesNoEmail%1,                Cadwallender%Bruce A.%coreGivesNoEmail%1,                Cinti%Sandro K.%coreGivesNoEmail%1,                Collins%Curtis D.%coreGivesNoEmail%1,                Goldberg%Janet%coreGivesNoEmail%1,                Holmes%Jennifer G.%coreGivesNoEmail%1,                Kim%Christopher%coreGivesNoEmail%1,                Krupansky%Frank%coreGivesNoEmail%1,                Lozon%Marie M.%coreGivesNoEmail%1,                Rodgers%Phillip E.%coreGivesNoEmail%1,                Shlafer%Jean%coreGivesNoEmail%1,                Wagner%Deborah%coreGivesNoEmail%1,                Wilkerson%William M.%coreGivesNoEmail%1,                Wright%Carrie M.%coreGivesNoEmail%1]'

# E4: update Authors value (cleaned data - extra whitespace normalization)
$ws.Range("E4").Value = '[Cathy%Campbell%xref no email%1,          Marianne%Baernholdt%xref no email%1]'
